# Edit 1: split the "ssh username@..." run into 3 runs, changing username -> grovesti,
# scoped to the specific paragraph (list item under "Per fare tunnel") to avoid touching
# the unrelated "username" occurrences elsewhere in the document.
$d = $word.ActiveDocument

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq "ssh username@sshpaolotti.studenti.math.unipd.it -L8080:tecweb:80 -L8022:tecweb:22") {
        $targetPara = $cand
        break
    }
}

if ($targetPara -ne $null) {
    $sshXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">ssh </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>grovesti</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>@sshpaolotti.studenti.math.unipd.it -L8080:tecweb:80 -L8022:tecweb:22</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $targetPara.Range.InsertXML($sshXml)
}

# Edit 2: mark the run that hosts the Filezilla screenshot drawing as NoProof (<w:noProof/>).
$picPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.InlineShapes.Count -gt 0) {
        $picPara = $cand
        break
    }
}

if ($picPara -ne $null) {
    $picPara.Range.NoProofing = 1

    # Edit 3: after the screenshot paragraph, append a page break, a short instruction
    # paragraph, the squadra.php link (split across several runs), and a trailing empty
    # paragraph.
    $picPara.Range.InsertParagraphAfter()
    $afterPicIndex = $picPara.Index + 1
    $newPara = $d.Paragraphs.Item($afterPicIndex)
    $tailXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Dall’interno trascinare i file da sinistra a destra nello screen sopra e poi collegarsi ade esempio su:</w:t></w:r></w:p><w:p><w:r><w:t>tecweb.studenti.math.unipd.it/grovest</w:t></w:r><w:r><w:t>i</w:t></w:r><w:r><w:t>/sq</w:t></w:r><w:r><w:t>uadra.php</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($tailXml)
}
